$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells that look numeric stay stored as text,
# matching the original inline-string formatting (e.g. "1.000", "51.30").
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value2 = '23.961.02'
$ws.Range("E2").Value = '  +0.51%  '
$ws.Range("D3").Value2 = '1.657.35'
$ws.Range("E3").Value = '  +2.17%  '
$ws.Range("D4").Value2 = '1.003'
$ws.Range("E4").Value = '  -0.42%  '
$ws.Range("D5").Value2 = '309.75'
$ws.Range("E5").Value = '  +1.13%  '
$ws.Range("E6").Value = '  -0.42%  '
$ws.Range("D7").Value2 = '0.3899'
$ws.Range("E7").Value = '  -0.28%  '
$ws.Range("D8").Value2 = '0.3867'
$ws.Range("E8").Value = '  +1.24%  '
$ws.Range("D9").Value2 = '51.30'
$ws.Range("E9").Value = '  +3.20%  '
$ws.Range("D10").Value2 = '1.362'
$ws.Range("E10").Value = '  +0.34%  '
$ws.Range("D11").Value2 = '1.003'
$ws.Range("E11").Value = '  -0.19%  '
$ws.Range("D12").Value2 = '0.08489'
$ws.Range("E12").Value = '  +1.16%  '
$ws.Range("D13").Value2 = '23.87'
$ws.Range("E13").Value = '  +0.63%  '
$ws.Range("D14").Value2 = '7.204'
$ws.Range("E14").Value = '  +3.55%  '
$ws.Range("D15").Value2 = '8.011'
$ws.Range("E15").Value = '  +8.22%  '
$ws.Range("D16").Value2 = '0.00001311'
$ws.Range("E16").Value = '  +3.71%  '
$ws.Range("D17").Value2 = '1.657.02'
$ws.Range("E17").Value = '  +2.43%  '
$ws.Range("D18").Value2 = '94.52'
$ws.Range("E18").Value = '  +2.13%  '
$ws.Range("D19").Value2 = '0.06985'
$ws.Range("E19").Value = '  +1.18%  '
$ws.Range("D20").Value2 = '19.90'
$ws.Range("E20").Value = '  +0.56%  '
$ws.Range("D21").Value2 = '6.971'
$ws.Range("E21").Value = '  +2.33%  '
$ws.Range("E22").Value = '  -0.50%  '
$ws.Range("D23").Value2 = '13.64'
$ws.Range("E23").Value = '  +2.40%  '
$ws.Range("D24").Value2 = '23.968.88'
$ws.Range("E24").Value = '  +0.45%  '
$ws.Range("D25").Value2 = '2.489'
$ws.Range("E25").Value = '  +4.38%  '
$ws.Range("D26").Value2 = '3.082'
$ws.Range("E26").Value = '  +8.06%  '
$ws.Range("D27").Value2 = '22.26'
$ws.Range("E27").Value = '  +1.20%  '
$ws.Range("D28").Value2 = '153.25'
$ws.Range("E28").Value = '  -2.81%  '
$ws.Range("E29").Value = '  +1.08%  '
$ws.Range("D30").Value2 = '5.312'
$ws.Range("E30").Value = '  +1.68%  '
$ws.Range("D31").Value2 = '7.943'
$ws.Range("E31").Value = '  +4.09%  '
$ws.Range("E32").Value = '  +0.94%  '
$ws.Range("D33").Value2 = '1.846.38'
$ws.Range("E33").Value = '  +2.23%  '
$ws.Range("D34").Value2 = '1.040'
$ws.Range("E34").Value = '  +9.34%  '
$ws.Range("D35").Value2 = '0.08135'
$ws.Range("E35").Value = '  +2.92%  '
$ws.Range("D36").Value2 = '0.02993'
$ws.Range("D37").Value2 = '11.07'
$ws.Range("E37").Value = '  +7.53%  '
$ws.Range("D38").Value2 = '6.709'
$ws.Range("E38").Value = '  +2.60%  '
$ws.Range("D39").Value2 = '0.2697'
$ws.Range("E39").Value = '  +1.98%  '
$ws.Range("D40").Value2 = '0.09149'
$ws.Range("E40").Value = '  +0.51%  '
$ws.Range("D41").Value2 = '13.64'
$ws.Range("E41").Value = '  +1.60%  '
$ws.Range("D42").Value2 = '0.7543'
$ws.Range("E42").Value = '  +1.75%  '
$ws.Range("D43").Value2 = '1.419'
$ws.Range("E43").Value = '  +0.27%  '
$ws.Range("D44").Value2 = '16.52'
$ws.Range("E44").Value = '  +3.23%  '
$ws.Range("D45").Value2 = '0.7008'
$ws.Range("E45").Value = '  +3.11%  '
$ws.Range("D46").Value2 = '2.484'
$ws.Range("E46").Value = '  +2.37%  '
$ws.Range("D47").Value2 = '4.087'
$ws.Range("E47").Value = '  +0.87%  '
$ws.Range("D48").Value2 = '1.000'
$ws.Range("E48").Value = '  -0.50%  '
$ws.Range("D49").Value2 = '0.08286'
$ws.Range("E49").Value = '  +0.70%  '
$ws.Range("D50").Value2 = '135.40'
$ws.Range("E50").Value = '  +2.45%  '
$ws.Range("D51").Value2 = '1.233'
$ws.Range("E51").Value = '  -0.70%  '

# Restore default (un-styled) formatting on the price column so the
# cells keep no explicit style index, same as before the edit.
$priceRange.Style = "Normal"

